$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Cells.Item(6, 1).Value = 9
$ws.Cells.Item(6, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(6, 3).Value = "Metropolitana"
$ws.Cells.Item(6, 4).Value = 44799
$ws.Cells.Item(6, 5).Value = 13
$ws.Cells.Item(6, 6).Value = 100112022
$ws.Cells.Item(6, 7).Value = "Arveja Verde"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 20
$ws.Cells.Item(6, 11).Value = 41000
$ws.Cells.Item(6, 12).Value = 41000
$ws.Cells.Item(6, 13).Value = 41000
$ws.Cells.Item(6, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(6, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(6, 16).Value = 1640
$ws.Cells.Item(6, 17).Value = 25
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# Row 7
$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(7, 3).Value = "Metropolitana"
$ws.Cells.Item(7, 4).Value = 44208
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 100112022
$ws.Cells.Item(7, 7).Value = "Arveja Verde"
$ws.Cells.Item(7, 8).Value = "Perfection"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 30
$ws.Cells.Item(7, 11).Value = 27000
$ws.Cells.Item(7, 12).Value = 27000
$ws.Cells.Item(7, 13).Value = 27000
$ws.Cells.Item(7, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Carahue"
$ws.Cells.Item(7, 16).Value = 1080
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Row 8
$ws.Cells.Item(8, 1).Value = 9
$ws.Cells.Item(8, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(8, 3).Value = "Metropolitana"
$ws.Cells.Item(8, 4).Value = 44201
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(8, 6).Value = 100112022
$ws.Cells.Item(8, 7).Value = "Arveja Verde"
$ws.Cells.Item(8, 8).Value = "Perfection"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 30
$ws.Cells.Item(8, 11).Value = 28000
$ws.Cells.Item(8, 12).Value = 28000
$ws.Cells.Item(8, 13).Value = 28000
$ws.Cells.Item(8, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(8, 15).Value = "Carahue"
$ws.Cells.Item(8, 16).Value = 1120
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Row 9
$ws.Cells.Item(9, 1).Value = 9
$ws.Cells.Item(9, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(9, 3).Value = "Metropolitana"
$ws.Cells.Item(9, 4).Value = 44536
$ws.Cells.Item(9, 5).Value = 13
$ws.Cells.Item(9, 6).Value = 100112022
$ws.Cells.Item(9, 7).Value = "Arveja Verde"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 43
$ws.Cells.Item(9, 11).Value = 16000
$ws.Cells.Item(9, 12).Value = 17000
$ws.Cells.Item(9, 13).Value = 16512
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Región del Maule"
$ws.Cells.Item(9, 16).Value = 660
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(10, 3).Value = "Metropolitana"
$ws.Cells.Item(10, 4).Value = 44537
$ws.Cells.Item(10, 5).Value = 13
$ws.Cells.Item(10, 6).Value = 100112022
$ws.Cells.Item(10, 7).Value = "Arveja Verde"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 61
$ws.Cells.Item(10, 11).Value = 15000
$ws.Cells.Item(10, 12).Value = 17000
$ws.Cells.Item(10, 13).Value = 15984
$ws.Cells.Item(10, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(10, 15).Value = "Región del Maule"
$ws.Cells.Item(10, 16).Value = 639
$ws.Cells.Item(10, 17).Value = 25
$ws.Cells.Item(10, 18).Value = "Hortaliza"

# Row 11
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(11, 3).Value = "Metropolitana"
$ws.Cells.Item(11, 4).Value = 44301
$ws.Cells.Item(11, 5).Value = 13
$ws.Cells.Item(11, 6).Value = 100112022
$ws.Cells.Item(11, 7).Value = "Arveja Verde"
$ws.Cells.Item(11, 8).Value = "Perfection"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 32000
$ws.Cells.Item(11, 12).Value = 32000
$ws.Cells.Item(11, 13).Value = 32000
$ws.Cells.Item(11, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(11, 16).Value = 1280
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Row 12
$ws.Cells.Item(12, 1).Value = 9
$ws.Cells.Item(12, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(12, 3).Value = "Metropolitana"
$ws.Cells.Item(12, 4).Value = 44519
$ws.Cells.Item(12, 5).Value = 13
$ws.Cells.Item(12, 6).Value = 100112022
$ws.Cells.Item(12, 7).Value = "Arveja Verde"
$ws.Cells.Item(12, 8).Value = "Sin especificar"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 34
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15500
$ws.Cells.Item(12, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(12, 15).Value = "Región Metropolitana"
$ws.Cells.Item(12, 16).Value = 620
$ws.Cells.Item(12, 17).Value = 25
$ws.Cells.Item(12, 18).Value = "Hortaliza"

# Row 13
$ws.Cells.Item(13, 1).Value = 9
$ws.Cells.Item(13, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44425
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = 100112022
$ws.Cells.Item(13, 7).Value = "Arveja Verde"
$ws.Cells.Item(13, 8).Value = "Perfection"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 16
$ws.Cells.Item(13, 11).Value = 35000
$ws.Cells.Item(13, 12).Value = 36000
$ws.Cells.Item(13, 13).Value = 35500
$ws.Cells.Item(13, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(13, 16).Value = 1420
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"

# Row 14
$ws.Cells.Item(14, 1).Value = 9
$ws.Cells.Item(14, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(14, 3).Value = "Metropolitana"
$ws.Cells.Item(14, 4).Value = 44469
$ws.Cells.Item(14, 5).Value = 13
$ws.Cells.Item(14, 6).Value = 100112022
$ws.Cells.Item(14, 7).Value = "Arveja Verde"
$ws.Cells.Item(14, 8).Value = "Perfection"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 16
$ws.Cells.Item(14, 11).Value = 25000
$ws.Cells.Item(14, 12).Value = 26000
$ws.Cells.Item(14, 13).Value = 25500
$ws.Cells.Item(14, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(14, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(14, 16).Value = 1020
$ws.Cells.Item(14, 17).Value = 25
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# Row 15
$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44452
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = 100112022
$ws.Cells.Item(15, 7).Value = "Arveja Verde"
$ws.Cells.Item(15, 8).Value = "Perfection"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 25
$ws.Cells.Item(15, 11).Value = 37000
$ws.Cells.Item(15, 12).Value = 38000
$ws.Cells.Item(15, 13).Value = 37480
$ws.Cells.Item(15, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(15, 16).Value = 1499
$ws.Cells.Item(15, 17).Value = 25
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Row 16
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44643
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = 100112022
$ws.Cells.Item(16, 7).Value = "Arveja Verde"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 20
$ws.Cells.Item(16, 11).Value = 31000
$ws.Cells.Item(16, 12).Value = 32000
$ws.Cells.Item(16, 13).Value = 31500
$ws.Cells.Item(16, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(16, 15).Value = "Carahue"
$ws.Cells.Item(16, 16).Value = 1260
$ws.Cells.Item(16, 17).Value = 25
$ws.Cells.Item(16, 18).Value = "Hortaliza"

# Row 17
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44671
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100112022
$ws.Cells.Item(17, 7).Value = "Arveja Verde"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 25
$ws.Cells.Item(17, 11).Value = 26000
$ws.Cells.Item(17, 12).Value = 27000
$ws.Cells.Item(17, 13).Value = 26520
$ws.Cells.Item(17, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(17, 15).Value = "Carahue"
$ws.Cells.Item(17, 16).Value = 1061
$ws.Cells.Item(17, 17).Value = 25
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Row 18
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44461
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112022
$ws.Cells.Item(18, 7).Value = "Arveja Verde"
$ws.Cells.Item(18, 8).Value = "Perfection"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 16
$ws.Cells.Item(18, 11).Value = 28000
$ws.Cells.Item(18, 12).Value = 29000
$ws.Cells.Item(18, 13).Value = 28500
$ws.Cells.Item(18, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(18, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(18, 16).Value = 1140
$ws.Cells.Item(18, 17).Value = 25
$ws.Cells.Item(18, 18).Value = "Hortaliza"

# Row 19
$ws.Cells.Item(19, 1).Value = 9
$ws.Cells.Item(19, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = 44229
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = 100112022
$ws.Cells.Item(19, 7).Value = "Arveja Verde"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 24000
$ws.Cells.Item(19, 12).Value = 24000
$ws.Cells.Item(19, 13).Value = 24000
$ws.Cells.Item(19, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value = "Carahue"
$ws.Cells.Item(19, 16).Value = 960
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"

# Row 20
$ws.Cells.Item(20, 1).Value = 9
$ws.Cells.Item(20, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(20, 3).Value = "Metropolitana"
$ws.Cells.Item(20, 4).Value = 44222
$ws.Cells.Item(20, 5).Value = 13
$ws.Cells.Item(20, 6).Value = 100112022
$ws.Cells.Item(20, 7).Value = "Arveja Verde"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 24000
$ws.Cells.Item(20, 12).Value = 24000
$ws.Cells.Item(20, 13).Value = 24000
$ws.Cells.Item(20, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Carahue"
$ws.Cells.Item(20, 16).Value = 960
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"

# Row 21
$ws.Cells.Item(21, 1).Value = 9
$ws.Cells.Item(21, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(21, 3).Value = "Metropolitana"
$ws.Cells.Item(21, 4).Value = 44475
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = 100112022
$ws.Cells.Item(21, 7).Value = "Arveja Verde"
$ws.Cells.Item(21, 8).Value = "Perfection"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 16
$ws.Cells.Item(21, 11).Value = 24000
$ws.Cells.Item(21, 12).Value = 25000
$ws.Cells.Item(21, 13).Value = 24500
$ws.Cells.Item(21, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(21, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(21, 16).Value = 980
$ws.Cells.Item(21, 17).Value = 25
$ws.Cells.Item(21, 18).Value = "Hortaliza"

# Row 22
$ws.Cells.Item(22, 1).Value = 9
$ws.Cells.Item(22, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(22, 3).Value = "Metropolitana"
$ws.Cells.Item(22, 4).Value = 44831
$ws.Cells.Item(22, 5).Value = 13
$ws.Cells.Item(22, 6).Value = 100112022
$ws.Cells.Item(22, 7).Value = "Arveja Verde"
$ws.Cells.Item(22, 8).Value = "Perfection"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 25
$ws.Cells.Item(22, 11).Value = 30000
$ws.Cells.Item(22, 12).Value = 30000
$ws.Cells.Item(22, 13).Value = 30000
$ws.Cells.Item(22, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(22, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(22, 16).Value = 1200
$ws.Cells.Item(22, 17).Value = 25
$ws.Cells.Item(22, 18).Value = "Hortaliza"

# Row 140
$ws.Cells.Item(140, 1).Value = 9
$ws.Cells.Item(140, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(140, 3).Value = "Metropolitana"
$ws.Cells.Item(140, 4).Value = 44832
$ws.Cells.Item(140, 5).Value = 13
$ws.Cells.Item(140, 6).Value = 100112022
$ws.Cells.Item(140, 7).Value = "Arveja Verde"
$ws.Cells.Item(140, 8).Value = "Perfection"
$ws.Cells.Item(140, 9).Value = "Primera"
$ws.Cells.Item(140, 10).Value = 25
$ws.Cells.Item(140, 11).Value = 30000
$ws.Cells.Item(140, 12).Value = 30000
$ws.Cells.Item(140, 13).Value = 30000
$ws.Cells.Item(140, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(140, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(140, 16).Value = 1200
$ws.Cells.Item(140, 17).Value = 25
$ws.Cells.Item(140, 18).Value = "Hortaliza"

# Ensure new row 140 date cell uses the same date style as other Fecha column cells
$ws.Cells.Item(140, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
